$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original column B data values (rows 2-6) before the shift, since
# the newly-inserted column duplicates column B's data (the new outcome
# column was added right after the existing "Programming I" skill column).
$origB2 = $ws.Range("B2").Value2
$origB3 = $ws.Range("B3").Value2
$origB4 = $ws.Range("B4").Value2
$origB5 = $ws.Range("B5").Value2
$origB6 = $ws.Range("B6").Value2

# Insert a new column before column B, shifting B:G -> C:H.
$ws.Range("B1").EntireColumn.Insert()

# New header text for the freshly inserted column B1.
$ws.Range("B1").Value = "Understand and implement advanced programming concepts building upon the foundational knowledge acquired in Programming I."

# Match the header style used by the rest of row 1 (bold, bordered, centered).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The new column B (rows 2-6) duplicates the original column B values.
$ws.Range("B2").Value = $origB2
$ws.Range("B3").Value = $origB3
$ws.Range("B4").Value = $origB4
$ws.Range("B5").Value = $origB5
$ws.Range("B6").Value = $origB6

# The insert copied column A's bordered style onto the new column's data
# rows; the data rows should remain unstyled, matching the other columns.
$ws.Range("B2:B6").ClearFormats()
